# Update cryptos list cell values to match the latest scraped data
# (commit: "Updated cryptos list on Sun Nov  5 05:30:35 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.613.35'
$ws.Range('E2').Value = '  +2.06%  '
$ws.Range('D3').Value = '1.908.48'
$ws.Range('E3').Value = '  +3.71%  '
$ws.Range('E4').Value = '  +0.55%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.39'
$ws.Range('E5').Value = '  +5.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.634'
$ws.Range('E6').Value = '  +2.67%  '
$ws.Range('E7').Value = '  +0.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.50'
$ws.Range('E8').Value = '  +3.65%  '
$ws.Range('E9').Value = '  +3.59%  '
$ws.Range('E10').Value = '  +2.52%  '
$ws.Range('E11').Value = '  +1.34%  '
$ws.Range('D12').Value = '2.184.19'
$ws.Range('E12').Value = '  +3.60%  '
$ws.Range('E13').Value = '  +10.04%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.915.56'
$ws.Range('E14').Value = '  +4.04%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.694'
$ws.Range('E15').Value = '  +3.52%  '
$ws.Range('D17').Value = '35.594.13'
$ws.Range('E17').Value = '  +1.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '72.21'
$ws.Range('E18').Value = '  +3.40%  '
$ws.Range('D19').Value = '0.0₃0811'
$ws.Range('E19').Value = '  +2.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '244.48'
$ws.Range('E20').Value = '  +1.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.54'
$ws.Range('E21').Value = '  +3.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.92'
$ws.Range('E22').Value = '  +3.70%  '
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('E24').Value = '  +1.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.50'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('E26').Value = '  +29.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.49'
$ws.Range('E27').Value = '  +8.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.99'
$ws.Range('E28').Value = '  +3.44%  '
$ws.Range('E29').Value = '  +1.91%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.10'
$ws.Range('E30').Value = '  +3.99%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0566'
$ws.Range('E31').Value = '  +2.92%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.948'
$ws.Range('E32').Value = '  +27.60%  '
$ws.Range('E33').Value = '  +0.48%  '
$ws.Range('E34').Value = '  +6.35%  '
$ws.Range('E35').Value = '  +7.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.04'
$ws.Range('E36').Value = '  +5.23%  '
$ws.Range('E37').Value = '  +5.78%  '
$ws.Range('E38').Value = '  +4.65%  '
$ws.Range('E39').Value = '  +4.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '91.61'
$ws.Range('E40').Value = '  +2.08%  '
$ws.Range('D41').Value = '1.359.90'
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.27'
$ws.Range('E42').Value = '  +5.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0597'
$ws.Range('E43').Value = '  +12.80%  '
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '48.16'
$ws.Range('E44').Value = '  +42.19%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.36'
$ws.Range('E45').Value = '  +4.68%  '
$ws.Range('B46').Value = 'Gas'
$ws.Range('C46').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.91'
$ws.Range('E46').Value = '  +16.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.43'
$ws.Range('E47').Value = '  +0.94%  '
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('E49').Value = '  +5.81%  '
$ws.Range('D50').Value = '2.093.67'
$ws.Range('E50').Value = '  +3.37%  '
$ws.Range('E51').Value = '  +4.35%  '
